$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 336; $r -le 442; $r++) {
    $ws.Cells.Item($r, 2).Value = "y"
}

$ws.Cells.Item(441, 1).Copy()
$ws.Cells.Item(442, 1).PasteSpecial(-4122)

$win = $excel.ActiveWindow
$win.ScrollRow = 399
$win.ScrollColumn = 1
$ws.Range("B435").Select()
